# Replace the Manager column (C) employee numbers with the manager's name,
# via a self-join on EMP_Number -> Name, then drop the row whose manager
# is NULL (KING, row 10) - mirrors an inner-join SQL query result.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Manager id -> manager name lookup, built from the employee table itself.
$managerNames = @{
    "7369" = "SMITH"
    "7499" = "ALLEN"
    "7521" = "WARD"
    "7566" = "JONES"
    "7654" = "MARTIN"
    "7698" = "BLAKE"
    "7782" = "CLARK"
    "7788" = "SCOTT"
    "7839" = "KING"
    "7844" = "TURNER"
    "7876" = "ADAMS"
    "7900" = "JAMES"
    "7902" = "FORD"
    "7934" = "MILLER"
}

$ws.Range("C2").Value = $managerNames["7902"]
$ws.Range("C3").Value = $managerNames["7698"]
$ws.Range("C4").Value = $managerNames["7698"]
$ws.Range("C5").Value = $managerNames["7839"]
$ws.Range("C6").Value = $managerNames["7698"]
$ws.Range("C7").Value = $managerNames["7839"]
$ws.Range("C8").Value = $managerNames["7839"]
$ws.Range("C9").Value = $managerNames["7566"]
# Row 10 (KING) has no manager - will be removed below.
$ws.Range("C11").Value = $managerNames["7698"]
$ws.Range("C12").Value = $managerNames["7788"]
$ws.Range("C13").Value = $managerNames["7698"]
$ws.Range("C14").Value = $managerNames["7566"]
$ws.Range("C15").Value = $managerNames["7782"]

# Drop the row with no manager (KING) - rows below shift up.
$ws.Rows(10).Delete()
